$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A1 value from 1 to 2
$ws.Range("A1").Value = 2

# Set the active selection to E3
$ws.Range("E3").Select()
